$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, B, C, D, E
$rows = @(
  @(2,  4392.86,   133.12,   4901.9,    148.54),
  @(3,  4837.83,   146.6,    606.73,    18.39),
  @(4,  444.97,    13.48,    -4295.17,  -130.15),
  @(5,  4435.8,    134.42,   5120.98,   155.18),
  @(6,  899.52,    27.26,    1131.7,    34.29),
  @(7,  -3536.28,  -107.16,  -3989.28,  -120.89),
  @(8,  4817.45,   145.98,   4808.9,    145.72),
  @(9,  2928.13,   88.73,    2889.89,   87.56999999999999),
  @(10, -1889.32,  -57.25,   -1919.01,  -58.15),
  @(11, 1373,      41.61,    2614.63,   79.23),
  @(12, 2753.59,   83.44,    1801.63,   54.59),
  @(13, 1380.59,   41.83,    -813,      -24.64),
  @(14, 808.87,    24.51,    935.24,    28.34),
  @(15, 817.63,    24.78,    810.4,     24.56),
  @(16, 8.76,      0.27,     -124.84,   -3.78)
)

foreach ($r in $rows) {
  $row = $r[0]
  $ws.Cells.Item($row, 2).Value = $r[1]
  $ws.Cells.Item($row, 3).Value = $r[2]
  $ws.Cells.Item($row, 4).Value = $r[3]
  $ws.Cells.Item($row, 5).Value = $r[4]
}
